$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.731.71"
$ws.Range("E2").Value = "  +0.15%  "
$ws.Range("D3").Value = "3.489.15"
$ws.Range("E3").Value = "  -0.13%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "592.55"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.33%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "171.58"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.68%  "
$ws.Range("E8").Value = "  -1.09%  "
$ws.Range("E9").Value = "  +2.67%  "
$ws.Range("E10").Value = "  -1.30%  "
$ws.Range("E11").Value = "  -2.02%  "
$ws.Range("D12").Value = "4.091.28"
$ws.Range("E12").Value = "  -0.18%  "
$ws.Range("E13").Value = "  -0.45%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "28.79"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.08%  "
$ws.Range("D15").Value = "66.764.50"
$ws.Range("E15").Value = "  +0.21%  "
$ws.Range("E16").Value = "  -0.79%  "
$ws.Range("D17").Value = "3.514.87"
$ws.Range("E17").Value = "  +0.66%  "
$ws.Range("E18").Value = "  -0.83%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.00"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.99%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "392.49"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "7.91"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.77%  "
$ws.Range("E22").Value = "  -0.88%  "
$ws.Range("E23").Value = "  +0.25%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.534"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.53%  "
$ws.Range("E25").Value = "  -1.82%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.15"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.29%  "
$ws.Range("E27").Value = "  -0.81%  "
$ws.Range("E28").Value = "  -0.06%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.17"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.88%  "
$ws.Range("E30").Value = "  -4.06%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "23.68"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.27%  "
$ws.Range("E33").Value = "  -1.39%  "
$ws.Range("E34").Value = "  -0.98%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "163.17"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.48%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.876"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.93%  "
$ws.Range("E37").Value = "  -1.87%  "
$ws.Range("E38").Value = "  +1.66%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "4.64"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.23%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0738"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.25%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "27.08"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.15%  "
$ws.Range("E42").Value = "  -1.28%  "
$ws.Range("D43").Value = "2.799.35"
$ws.Range("E43").Value = "  +0.57%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "42.71"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.18%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.54"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.36%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0301"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.73%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "336.02"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -4.58%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "34.60"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.95%  "
$ws.Range("E49").Value = "  -2.30%  "
$ws.Range("E50").Value = "  -1.45%  "
$ws.Range("E51").Value = "  -2.23%  "
